# Add a new "Диагноз" (Diagnosis) column to the patient table and drop the
# old helper COUNTIF formula that used to live in column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + data for column D.
$ws.Range("D1").Value = "Диагноз"
$ws.Range("D2").Value = "Обследование"

# The old formula cell (E1: =COUNTIF(A1:A10000,"<>")) is no longer needed.
$ws.Range("E1").ClearContents()

# Match the saved selection state from the edited workbook.
$ws.Range("D11").Select()
